$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.86756867170333862
$ws.Range("B1").Value = 2.219003438949585
$ws.Range("C1").Value = -1.218452930450439
$ws.Range("D1").Value = -4.0050206184387207
$ws.Range("E1").Value = -2.136902809143066
$ws.Range("F1").Value = -655.08083265941673

$ws.Range("A3").Value = -0.0476524047553539
$ws.Range("B3").Value = 1.741524815559387
$ws.Range("C3").Value = 1.9635529518127439
$ws.Range("D3").Value = -0.72651195526123047
$ws.Range("E3").Value = 2.9309113025665279
$ws.Range("F3").Value = 835.57526296724336

$ws.Range("A4").Value = -4.0830531120300293
$ws.Range("B4").Value = 0.70842981338500977
$ws.Range("C4").Value = -2.4846889972686772
$ws.Range("D4").Value = -3.2637290954589839
$ws.Range("E4").Value = -9.123042106628418
$ws.Range("F4").Value = -2680.348347935208

$ws.Range("A5").Value = -7.8781418800354004
$ws.Range("B5").Value = 1.4281715154647829
$ws.Range("C5").Value = -6.629636287689209
$ws.Range("D5").Value = -1.0083484649658201
$ws.Range("E5").Value = -14.087953567504879
$ws.Range("F5").Value = -4890.2353139030411

$ws.Range("A6").Value = 1.820960640907288
$ws.Range("B6").Value = 3.261458158493042
$ws.Range("C6").Value = 3.1642336845397949
$ws.Range("D6").Value = -10.9933967590332
$ws.Range("E6").Value = -2.7467401027679439
$ws.Range("F6").Value = -581.57982457169976

$ws.Range("A7").Value = 11.36473178863525
$ws.Range("B7").Value = 4.3121323585510254
$ws.Range("C7").Value = 18.505096435546879
$ws.Range("D7").Value = 4.9041824340820313
$ws.Range("E7").Value = 39.086143493652337
$ws.Range("F7").Value = 10781.242369742689

$ws.Range("A8").Value = -8.6794157028198242
$ws.Range("B8").Value = -2.90678882598877
$ws.Range("C8").Value = -5.6758871078491211
$ws.Range("D8").Value = -1.217256546020508
$ws.Range("E8").Value = -18.47934722900391
$ws.Range("F8").Value = -6138.6469825463982

$ws.Range("A9").Value = -17.688358306884769
$ws.Range("B9").Value = 0.032858673483133302
$ws.Range("C9").Value = -0.59555220603942871
$ws.Range("D9").Value = 34.269393920898438
$ws.Range("E9").Value = 16.018362045288089
$ws.Range("F9").Value = 5807.2170754401277

$ws.Range("A10").Value = -4.282341480255127
$ws.Range("B10").Value = 6.067842960357666
$ws.Range("C10").Value = -0.26504403352737432
$ws.Range("D10").Value = -0.98349380493164063
$ws.Range("E10").Value = 0.5369640588760376
$ws.Range("F10").Value = 172.81659372533409

$ws.Range("A11").Value = -7.4785289764404297
$ws.Range("B11").Value = 2.0431375503540039
$ws.Range("C11").Value = -0.9100680947303772
$ws.Range("D11").Value = -5.8367652893066406
$ws.Range("E11").Value = -12.18222141265869
$ws.Range("F11").Value = -2560.2843480053721

$ws.Range("A12").Value = -3.7520606517791748
$ws.Range("B12").Value = 3.2247834205627441
$ws.Range("C12").Value = 9.0975990295410156
$ws.Range("D12").Value = -4.0643086433410636
$ws.Range("E12").Value = 4.5060186386108398
$ws.Range("F12").Value = 1034.5972982684721
